# Bulk-upload block-level template: replace sample Block-A..E rows with
# BLOCK1..BLOCK6 rows (duplicated manager contacts, same as the new
# "create association" UI), lower-case the header row, add two more
# duplicate rows, a data-validation dropdown for BlockType, and a
# bordered 3-row placeholder block below the data (same as bulk-upload
# template produced by the app).

function BGR($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header row - lower-cased field names
# ---------------------------------------------------------------
$ws.Range("A1").Value = "Sno"
$ws.Range("B1").Value = "blockname"
$ws.Range("C1").Value = "blocktype"
$ws.Range("D1").Value = "units"
$ws.Range("E1").Value = "managername"
$ws.Range("F1").Value = "managermobileno"
$ws.Range("G1").Value = "manageremailid"

# ---------------------------------------------------------------
# Data rows: 6 blocks, alternating SAM/sam007@gmail.com and
# JPHN/john@gmail.com manager contacts, 4 units each.
# ---------------------------------------------------------------
$rows = @(
    @{ n=1; block="BLOCK1"; type="Residential and Commercial"; units=4; mgr="SAM";   mail="sam007@gmail.com" },
    @{ n=2; block="BLOCK2"; type="Commercial";                  units=4; mgr="JPHN";  mail="john@gmail.com" },
    @{ n=3; block="BLOCK3"; type="Residential and Commercial"; units=4; mgr="SATYA"; mail="sam007@gmail.com" },
    @{ n=4; block="BLOCK4"; type="Residential";                 units=4; mgr="SILVER";mail="john@gmail.com" },
    @{ n=5; block="BLOCK5"; type="Commercial";                  units=4; mgr="ANU";   mail="sam007@gmail.com" },
    @{ n=6; block="BLOCK6"; type="Residential and Commercial"; units=4; mgr="DADY";  mail="john@gmail.com" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $row.n
    $ws.Cells.Item($r, 2).Value = $row.block
    $ws.Cells.Item($r, 3).Value = $row.type
    $ws.Cells.Item($r, 4).Value = $row.units
    $ws.Cells.Item($r, 5).Value = $row.mgr
    $ws.Cells.Item($r, 6).Value = 9886819118
    $ws.Cells.Item($r, 7).Value = $row.mail
}

# A-column banding style alternates between the two existing numeric
# styles already used in the sheet (s=3 / s=4) - odd data rows get the
# first banding style, even rows the second, matching rows 2..7.
$ws.Range("A2").Style = "Normal"
$aStyles = @(3, 4, 3, 4, 3, 4)
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    if ($aStyles[$i] -eq 3) {
        $ws.Cells.Item($r, 1).NumberFormat = $ws.Range("A2").NumberFormat
    }
}

# Re-apply the row-3 (s=4) banded style onto row 3/5/7 and the row-2
# (s=3) banded style onto row 4/6 by copying formats from the rows that
# already carry the right style.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4,A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A5,A7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Row heights: rows 4-7 now use the sheet default (15.15); rows 8-10
# (the new blank placeholder block) use 14.7, same as the old rows.
# ---------------------------------------------------------------
$ws.Rows.Item(4).RowHeight = 15.15
$ws.Rows.Item(5).RowHeight = 15.15
$ws.Rows.Item(6).RowHeight = 15.15
$ws.Rows.Item(7).RowHeight = 15.15
$ws.Rows.Item(8).RowHeight = 14.7
$ws.Rows.Item(9).RowHeight = 14.7
$ws.Rows.Item(10).RowHeight = 14.7

# ---------------------------------------------------------------
# BlockType dropdown validation for the whole editable column,
# including the blank placeholder rows.
# ---------------------------------------------------------------
$ws.Range("C2:C8").Validation.Delete() | Out-Null
$ws.Range("C2:C8").Validation.Add(3, 1, 1, """Residential,Commercial,Residential and Commercial""") | Out-Null
$ws.Range("C2:C8").Validation.IgnoreBlank = $true
$ws.Range("C2:C8").Validation.InCellDropdown = $true

# ---------------------------------------------------------------
# Hyperlinks: rebuild G2:G7 from scratch so the mailto targets line up
# with the (possibly duplicated) manager e-mail shown in each row.
# ---------------------------------------------------------------
$ws.Cells.Hyperlinks.Delete() | Out-Null
$mails = @("sam007@gmail.com","john@gmail.com","sam007@gmail.com","john@gmail.com","sam007@gmail.com","john@gmail.com")
for ($i = 0; $i -lt 6; $i++) {
    $r = $i + 2
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), "mailto:" + $mails[$i]) | Out-Null
}

# ---------------------------------------------------------------
# New blank placeholder block (rows 8-10, cols A:G) - light fill with a
# thin grey outline, same look as the header/data banding above.
# ---------------------------------------------------------------
$white = BGR 255 255 255
$top   = BGR 165 165 165
$side  = BGR 170 170 170

$placeholder = $ws.Range("A8:G10")
$placeholder.Interior.Color = $white

# Top row (row 8): top border across, left border on A8, right border on G8
$topRow = $ws.Range("A8:G8")
$topRow.Borders.Item(8).LineStyle = 1
$topRow.Borders.Item(8).Weight = 2
$topRow.Borders.Item(8).Color = $top
$ws.Range("A8").Borders.Item(7).LineStyle = 1
$ws.Range("A8").Borders.Item(7).Weight = 2
$ws.Range("A8").Borders.Item(7).Color = $side
$ws.Range("G8").Borders.Item(10).LineStyle = 1
$ws.Range("G8").Borders.Item(10).Weight = 2
$ws.Range("G8").Borders.Item(10).Color = $side

# Middle row (row 9): left border on A9, right border on G9 only
$ws.Range("A9").Borders.Item(7).LineStyle = 1
$ws.Range("A9").Borders.Item(7).Weight = 2
$ws.Range("A9").Borders.Item(7).Color = $side
$ws.Range("G9").Borders.Item(10).LineStyle = 1
$ws.Range("G9").Borders.Item(10).Weight = 2
$ws.Range("G9").Borders.Item(10).Color = $side

# Bottom row (row 10): bottom border across, left border on A10, right border on G10
$bottomRow = $ws.Range("A10:G10")
$bottomRow.Borders.Item(9).LineStyle = 1
$bottomRow.Borders.Item(9).Weight = 2
$bottomRow.Borders.Item(9).Color = $side
$ws.Range("A10").Borders.Item(7).LineStyle = 1
$ws.Range("A10").Borders.Item(7).Weight = 2
$ws.Range("A10").Borders.Item(7).Color = $side
$ws.Range("G10").Borders.Item(10).LineStyle = 1
$ws.Range("G10").Borders.Item(10).Weight = 2
$ws.Range("G10").Borders.Item(10).Color = $side

# ---------------------------------------------------------------
# Column widths - minor autosize drift from the real edit.
# ---------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 24.21875
$ws.Columns.Item(7).ColumnWidth = 25.6640625

# ---------------------------------------------------------------
# Selection cursor, matching the author's last active cell.
# ---------------------------------------------------------------
$ws.Range("B7").Select() | Out-Null
